$d = $word.ActiveDocument

# Locate the list-paragraph ending with "The design of output bandpass is "
# (this paragraph's identity/formatting is what survives the edit).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "The design of output bandpass is ") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The very next paragraph holds "May a bandp need ass near output" in its
    # own run/paragraph (its own numbering + different run formatting).
    # Remove that whole paragraph (text + its own paragraph mark) so it
    # collapses away, leaving $target's paragraph properties intact.
    $extra = $target.Next()
    $extra.Range.Delete()

    # Append the rest of the sentence right after the existing run's text
    # ("...bandpass is " already ends in a space, so this yields the exact
    # final wording) without disturbing the existing run's identity/formatting.
    $insertPoint = $d.Range($target.Range.End - 1, $target.Range.End - 1)
    $insertPoint.InsertAfter("to some difficult.")
}
